# Applies the BOQ bill update: revises line-item quantities/descriptions
# (rows 8-20), removes the old "Grand Total" row, and rebuilds the totals
# footer (Grand Total Rs. / Tender Premium / NET PAYABLE) with new amounts.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("A8").Value = 'P. point'
$ws.Range("C8").Value = '15'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3'
$ws.Range("E8").Value = 'Medium point (up to 6 mtr.)'
$ws.Range("F8").Value = '472'
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '7080.00'

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '4.0'
$ws.Range("E9").Value = 'P & F ISI marked (IS :3854) 16 amp. flush type non modular switch CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F9").Value = '50'
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '2850.00'

# Row 10
$ws.Range("C10").Value = '95'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '10.0'
$ws.Range("E10").Value = 'Providing and fixing of   power plug point with non modular accessories as per PWD specification for electrical Works with  Galvanized   box of 1.2 mm thick  with earth terminal with suitable size phenolic laminated sheet (IS : 2036 -  1995) cover including cost of 16 amp. Switch (IS :3854) and 3/6 pin 16 amp. socket outlet  making connection , testing , etc. as required. . For specification of  Wiring accessories refer Chapter  E - 07 related item &  For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F10").Value = '303'
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '28785.00'

# Row 11
$ws.Range("A11").Value = ""
$ws.Range("C11").Value = '60'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '11.0'
$ws.Range("E11").Value = 'S&F following sizes (dia.) of ISI marked virgin material MMS ( IS:9537 P - III ) PVC conduit along with  ISI marked (IS:3419-1988) accessories as required  in  recess  including  cutting the wall, covering conduit and making good the same as required. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F11").Value = '0'
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '0.00'

# Row 12
$ws.Range("A12").Value = 'Mtr.'
$ws.Range("C12").Value = '45'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '19'
$ws.Range("E12").Value = '2 x 2.5 sq. mm. + 1x1.5sqmm'
$ws.Range("F12").Value = '81'
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '3645.00'

# Row 13
$ws.Range("A13").Value = 'Mtr.'
$ws.Range("C13").Value = '86'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23'
$ws.Range("E13").Value = '8 SWG G.I. ( Hot Dipped  ) Wire '
$ws.Range("F13").Value = '20'
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '1720.00'

# Row 14
$ws.Range("A14").Value = 'Each'
$ws.Range("C14").Value = '70'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25'
$ws.Range("E14").Value = '1200 mm Sweep BEE 1 Star rated (service value >=4.0 to < 4.5 )'
$ws.Range("F14").Value = '1890'
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '132300.00'

# Row 15
$ws.Range("A15").Value = ""
$ws.Range("C15").Value = '36'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '16.0'
$ws.Range("E15").Value = 'Providing & Fixing of IP20 SMD Mid Power LED batten type integrated light fixture made from Powder coated Extruded aluminium  housing with in built driver  , System lumen efficacy ≥ 110 lm/Watt output, internal surge protection of 2.5 KV with Short & Open circuit protection ,THD < 10% , P. F.≥0.95, CRI >80 , life time of minimum  50000 Burning Hours with , 70% of intial Lumen maintaned till life ends  , CCT 3000°K / 4000°K  / 5700°K /6000°K/6500°K (As per ANSI Bin) , Maximum power consumption should not more than the specified rating and Fixture shall be of  BIS standard and  trade mark certificate ( T.C.). Manufactures Word Mark/ Name Engraved/ Embossing/ Screen printing on housing. OEM must have its own in house NABL lab setup for all testing facilities for LED fixtures. (LM79 & LM80) certificate / Report from OEM shall be submitted.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F15").Value = '0'
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '0.00'

# Row 16
$ws.Range("A16").Value = 'Each'
$ws.Range("C16").Value = '39'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '27'
$ws.Range("E16").Value = '1170mm(+/-10%) LED batten with min. lumen output 2200 lm'
$ws.Range("F16").Value = '492'
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '19188.00'

# Row 17
$ws.Range("C17").Value = '45'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '31'
$ws.Range("E17").Value = 'Double pole MCB(With B/C curve tripping Characteristics)'

# Row 18
$ws.Range("A18").Value = ""
$ws.Range("C18").Value = '53'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.0'
$ws.Range("E18").Value = 'Providing & Fixing of Recessed/surface mounting heavy duty horizontal type Double Door ( Metal / Glazed )Distribution board with Metal end box made out from Galvanized steel / CRCA sheet not less then 1.2 mm thick  conforming to IS-8623-1 & 3 /  IEC 61439- 1 & 3, powder painted complete with reversible door (for double door DB only )100 amp.  insulated copper bus bar/shorting link , copper neutral link, copper earth link , color coded interconnecting wire set  of suitable rating and din bar,masking sheet,  making internal DB  terminations with copper lugs, Ferrules,  detachable gland plate, including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F18").Value = '0'
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '0.00'

# Row 19
$ws.Range("C19").Value = '54'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '34'
$ws.Range("E19").Value = 'Metal door (single phase) IK-09 and IP-43 with Metal end box'

# Row 20
$ws.Range("A20").Value = ""
$ws.Range("C20").Value = '78'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '36'
$ws.Range("E20").Value = 'Total'

# Row 21 previously held the "Grand Total" line; its data cells are cleared
# (A21 stays an empty string).
$ws.Range("B21:I21").ClearContents()

# The old, mostly-empty row 22 is removed so the totals footer (old rows
# 23-25: Grand Total Rs. / Tender Premium / NET PAYABLE AMOUNT) shifts up
# to become rows 22-24.
$ws.Rows(22).Delete()

# Refresh the Grand Total Rs. (row 22) and NET PAYABLE AMOUNT Rs. (row 24)
# amounts to reflect the revised line items above; Tender Premium (row 23)
# is unchanged at 0.00.
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "195568.00"
$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = "195568.00"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "195568.00"
$ws.Range("H24").NumberFormat = "@"
$ws.Range("H24").Value = "195568.00"
